$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = 44175
$startRow = 5815

$districts = @(
    @('Bánovce nad Bebravou', 2),
    @('Banská Bystrica', 89),
    @('Bardejov', 12),
    @('Bratislava', 183),
    @('Brezno', 35),
    @('Bytča', 12),
    @('Čadca', 32),
    @('Detva', 14),
    @('Dolný Kubín', 2),
    @('Dunajská Streda', 11),
    @('Galanta', 38),
    @('Gelnica', 7),
    @('Hlohovec', 29),
    @('Humenné', 11),
    @('Ilava', 40),
    @('Kežmarok', 24),
    @('Komárno', 17),
    @('Košice', 62),
    @('Košice - okolie', 61),
    @('Krupina', 11),
    @('Kysucké Nové Mesto', 12),
    @('Levice', 26),
    @('Levoča', 14),
    @('Liptovský Mikuláš', 2),
    @('Lučenec', 19),
    @('Malacky', 52),
    @('Martin', 23),
    @('Medzilaborce', 34),
    @('Michalovce', 49),
    @('Myjava', 4),
    @('Námestovo', 4),
    @('Nitra', 217),
    @('Nové Mesto nad Váhom', 12),
    @('Nové Zámky', 60),
    @('Partizánske', 2),
    @('Pezinok', 25),
    @('Piešťany', 15),
    @('Poltár', 2),
    @('Poprad', 64),
    @('Považská Bystrica', 48),
    @('Prešov', 55),
    @('Prievidza', 7),
    @('Púchov', 35),
    @('Revúca', 6),
    @('Rimavská Sobota', 31),
    @('Rožňava', 7),
    @('Ružomberok', 10),
    @('Sabinov', 11),
    @('Senec', 26),
    @('Senica', 17),
    @('Skalica', 31),
    @('Snina', 11),
    @('Sobrance', 7),
    @('Spišská Nová Ves', 55),
    @('Stará Ľubovňa', 31),
    @('Stropkov', 17),
    @('Svidník', 11),
    @('Šaľa', 26),
    @('Topoľčany', 8),
    @('Trebišov', 35),
    @('Trenčín', 121),
    @('Trnava', 42),
    @('Turčianske Teplice', 6),
    @('Tvrdošín', 7),
    @('Veľký Krtíš', 33),
    @('Vranov nad Topľou', 42),
    @('Zlaté Moravce', 6),
    @('Zvolen', 30),
    @('Žarnovica', 7),
    @('Žiar nad Hronom', 3),
    @('Žilina', 56),
)

for ($i = 0; $i -lt $districts.Count; $i++) {
    $r = $startRow + $i
    $name = $districts[$i][0]
    $cnt = $districts[$i][1]
    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $cnt
}

